$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (Price D, Volume(1h) E) updates, per commit diff
$updates = @(
    @{ Row = 2; D = '27.486.76'; E = '  -2.40%  ' },
    @{ Row = 3; D = '1.748.17'; E = '  -2.69%  ' },
    @{ Row = 4; D = '1.002'; E = '  -0.04%  ' },
    @{ Row = 5; D = '323.74'; E = '  +0.01%  ' },
    @{ Row = 6; D = $null; E = '  -0.02%  ' },
    @{ Row = 7; D = '0.4463'; E = '  +4.01%  ' },
    @{ Row = 8; D = '0.3604'; E = '  -0.81%  ' },
    @{ Row = 9; D = '0.07464'; E = '  -1.12%  ' },
    @{ Row = 10; D = $null; E = '  -5.97%  ' },
    @{ Row = 11; D = '1.090'; E = '  -2.58%  ' },
    @{ Row = 12; D = '1.001'; E = '  +0.00%  ' },
    @{ Row = 13; D = '20.52'; E = '  -5.48%  ' },
    @{ Row = 14; D = '6.002'; E = '  -2.67%  ' },
    @{ Row = 15; D = '7.100'; E = '  -3.48%  ' },
    @{ Row = 16; D = '1.749.96'; E = '  -3.54%  ' },
    @{ Row = 17; D = '92.09'; E = '  -0.56%  ' },
    @{ Row = 18; D = '0.00001057'; E = '  -1.33%  ' },
    @{ Row = 19; D = '0.06403'; E = '  +0.92%  ' },
    @{ Row = 20; D = $null; E = '  +0.03%  ' },
    @{ Row = 21; D = '16.75'; E = '  -2.77%  ' },
    @{ Row = 22; D = $null; E = '  -2.53%  ' },
    @{ Row = 23; D = '27.531.88'; E = '  -2.27%  ' },
    @{ Row = 24; D = $null; E = '  -2.56%  ' },
    @{ Row = 25; D = '2.111'; E = '  -3.23%  ' },
    @{ Row = 26; D = '161.94'; E = '  +1.69%  ' },
    @{ Row = 27; D = '20.35'; E = '  -0.03%  ' },
    @{ Row = 28; D = '1.950.49'; E = '  -3.29%  ' },
    @{ Row = 29; D = '2.067'; E = '  -7.47%  ' },
    @{ Row = 30; D = '124.36'; E = '  -2.72%  ' },
    @{ Row = 31; D = '1.077'; E = '  -8.29%  ' },
    @{ Row = 32; D = '3.656'; E = '  +3.17%  ' },
    @{ Row = 33; D = '0.09015'; E = '  -0.08%  ' },
    @{ Row = 34; D = '5.475'; E = '  -6.48%  ' },
    @{ Row = 35; D = '11.94'; E = '  -6.29%  ' },
    @{ Row = 36; D = '0.02293'; E = '  -2.76%  ' },
    @{ Row = 37; D = '0.2078'; E = '  -1.98%  ' },
    @{ Row = 38; D = '0.6326'; E = '  -2.61%  ' },
    @{ Row = 39; D = '0.05964'; E = '  -2.56%  ' },
    @{ Row = 40; D = '4.909'; E = '  -3.83%  ' },
    @{ Row = 41; D = '1.204'; E = '  +0.95%  ' },
    @{ Row = 42; D = $null; E = '  +0.07%  ' },
    @{ Row = 43; D = '1.385'; E = '  -3.00%  ' },
    @{ Row = 44; D = '7.737'; E = '  -2.69%  ' },
    @{ Row = 45; D = '13.16'; E = '  -3.40%  ' },
    @{ Row = 46; D = '0.5863'; E = '  -2.76%  ' },
    @{ Row = 47; D = '3.698'; E = '  -0.29%  ' },
    @{ Row = 48; D = '120.93'; E = '  -3.61%  ' },
    @{ Row = 49; D = '1.940'; E = '  -3.10%  ' },
    @{ Row = 50; D = '1.147'; E = '  -0.47%  ' },
    @{ Row = 51; D = '0.06850'; E = '  -1.72%  ' }
)

foreach ($item in $updates) {
    $r = $item.Row
    if ($null -ne $item.D) {
        $dCell = $ws.Cells.Item($r, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $item.D
        $dCell.Style = "Normal"
    }
    $eCell = $ws.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $item.E
    $eCell.Style = "Normal"
}
